# Apply the documented changes:
# 1. Insert a new paragraph right after the title (Heading1) with a
#    "Meta description" lead-in (bold) followed by the descriptive text.
# 2. Remove the duplicate bold "Play Dead or Alive 2 Free..." paragraph
#    that used to sit near the end of the document (a second copy of the
#    title text that was re-used as a meta-description heading there).
# 3. Replace the text of the final (italic) paragraph with the new
#    image-generation prompt, keeping its italic formatting intact.

$d = $word.ActiveDocument

$titleText = "Play Dead or Alive 2 Free: Varied Gameplay and Beautiful Graphics!"
$oldMetaText = "Dead or Alive 2 is a Western-themed slot game with various entertaining modes and beautiful graphics. Play free and win big today!"
$newImagePrompt = "Please create a cartoon style image featuring a happy Maya warrior with glasses. The warrior should be holding a revolver in one hand and a wanted poster in the other, standing in front of a Wild West village inside a canyon, complete with cattle skulls, bottles, and carriage wheels. The background should have storm clouds approaching, and there should be creaking and barking dogs. The image should capture the excitement and adventure of the Wild West, as well as the fun and entertainment of playing a slot game."

# --- 1. Insert the "Meta description" paragraph right after the title ---
# Locate the title paragraph (the very first occurrence of the title text).
$titleRange = $d.Content
[void]$titleRange.Find.Execute($titleText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$titlePara = $titleRange.Paragraphs(1)

[void]$titlePara.Range.InsertParagraphAfter()
$metaPara = $titlePara.Next()

# Build the new paragraph's content exactly (empty run + bold lead-in run +
# plain-text run) via a raw OOXML fragment, so the resulting markup mirrors
# the structure used throughout the rest of the document.
$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: ' + $oldMetaText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$metaPara.Range.InsertXML($metaXml)

# --- 2. Delete the old duplicate bold title paragraph near the end ---
# Find the *second* occurrence of the title text (the first is the real
# Heading1 title untouched by this step) and remove its whole paragraph.
$dupRange = $d.Content
[void]$dupRange.Find.Execute($titleText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dupRange.Collapse(0)
[void]$dupRange.Find.Execute($titleText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dupPara = $dupRange.Paragraphs(1)
[void]$dupPara.Range.Delete()

# --- 3. Replace the text of the final italic paragraph with the new prompt ---
# Scope the search/replace to just the last paragraph so the earlier
# "Meta description" paragraph (which also contains the old description
# text) is left untouched.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$find = $lastPara.Range.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
[void]$find.Execute($oldMetaText, $true, $false, $false, $false, $false, $true, 1, $false, $newImagePrompt, 2)
